$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F5").Value = 69
$ws.Range("H5").Value = 69
$ws.Range("E7").Value = 27
$ws.Range("E9").Value = 8
$ws.Range("E10").Value = 380
$ws.Range("F10").Value = 174
$ws.Range("H10").Value = 174
$ws.Range("E11").Value = 249
$ws.Range("F11").Value = 134
$ws.Range("H11").Value = 134
$ws.Range("E12").Value = 366
$ws.Range("E14").Value = 94
$ws.Range("F14").Value = 46
$ws.Range("H14").Value = 46
$ws.Range("E15").Value = 122
$ws.Range("E16").Value = 162
$ws.Range("E17").Value = 70
$ws.Range("E20").Value = 74
$ws.Range("E21").Value = 119
$ws.Range("E22").Value = 142
$ws.Range("E23").Value = 164
$ws.Range("F23").Value = 72
$ws.Range("H23").Value = 72
$ws.Range("E24").Value = 164
$ws.Range("E25").Value = 192
$ws.Range("F25").Value = 88
$ws.Range("H25").Value = 88
$ws.Range("E26").Value = 112
$ws.Range("F26").Value = 66
$ws.Range("H26").Value = 66
$ws.Range("E27").Value = 253
$ws.Range("F27").Value = 121
$ws.Range("H27").Value = 121
$ws.Range("E28").Value = 149
$ws.Range("E29").Value = 139
$ws.Range("E30").Value = 164
$ws.Range("E32").Value = 152
$ws.Range("E33").Value = 233
$ws.Range("F33").Value = 113
$ws.Range("H33").Value = 113
$ws.Range("E34").Value = 169
$ws.Range("F34").Value = 99
$ws.Range("H34").Value = 99
$ws.Range("E35").Value = 111
$ws.Range("F35").Value = 70
$ws.Range("H35").Value = 70
$ws.Range("E37").Value = 122
$ws.Range("E39").Value = 153
$ws.Range("E40").Value = 210
$ws.Range("E41").Value = 306
$ws.Range("F41").Value = 130
$ws.Range("H41").Value = 130
$ws.Range("E42").Value = 272
$ws.Range("E43").Value = 93
$ws.Range("F43").Value = 45
$ws.Range("H43").Value = 45
$ws.Range("E45").Value = 109
$ws.Range("E46").Value = 238
$ws.Range("F46").Value = 130
$ws.Range("H46").Value = 130
$ws.Range("E47").Value = 349
$ws.Range("E48").Value = 161
$ws.Range("E49").Value = 232
$ws.Range("E50").Value = 203
$ws.Range("F50").Value = 71
$ws.Range("H50").Value = 71
$ws.Range("E51").Value = 187
$ws.Range("E52").Value = 23
